# remove preprocess image ocr
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,  "77A-247.01", "Top",    "2025-11-02 12:55:11"),
    @(4,  "77A-247.01", "Top",    "2025-11-02 14:11:38"),
    @(1,  "CN4",         "Bottom", "2025-11-02 14:21:41"),
    @(1,  "322",         "Bottom", "2025-11-02 14:25:49"),
    @(1,  "77C-226.75",  "Bottom", "2025-11-02 14:56:30"),
    @(35, "77A-247.01",  "Top",    "2025-11-02 14:56:45")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# "322" in B5 looks numeric; force it to stay text like the source OCR value.
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "322"
$ws.Cells.Item(5, 2).ClearFormats()
